$d = $word.ActiveDocument

# Locate the paragraph that ends with "Receiving the data on the ESP32 and Connect to the WIFI"
$rng = $d.Content
$found = $rng.Find.Execute("Receiving the data on the ESP32 and Connect to the WIFI", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $lines = @(
        "20/02/2021",
        "Storing value in flash memory in esp32",
        "21/02/2021",
        "Connecting the ESP32 to WIFI and server",
        "Fixing few bugs in the app",
        "Researching about the ESP32 sleep modes"
    )

    $insertPoint = $rng.End

    foreach ($line in $lines) {
        $d.Range($insertPoint, $insertPoint).InsertParagraphAfter()
        $newRange = $d.Range($insertPoint + 1, $insertPoint + 1)
        $newRange.Text = $line
        $insertPoint = $insertPoint + 1 + $line.Length
    }
}
